# Rename the "geo" axis/column header to "country" across the whole
# workbook. The header text "geo" appears in column A of every sheet
# ("pop", "births", "deaths", "pop_births_deaths", "pop_missing_axis_name",
# "pop_missing_values", "pop_narrow_format") -- including the repeated
# header rows inside "pop_births_deaths" (rows 1, 9 and 17). Replacing every
# occurrence updates the shared string used by all of them in one go.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rowCount = $used.Rows.Count
    $colCount = $used.Columns.Count
    for ($r = 1; $r -le $rowCount; $r++) {
        for ($c = 1; $c -le $colCount; $c++) {
            $cell = $used.Cells.Item($r, $c)
            if ($cell.Value2 -eq "geo") {
                $cell.Value2 = "country"
            }
        }
    }
}
